$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list: refresh Price (D) and Volume(1h) (E) columns.
# Values are set with a leading apostrophe to force text (matching the
# original inline-string / text-typed cells), then ClearFormats() strips
# the quote-prefix formatting Excel applies so the cell style is left
# untouched (same default/general style as before the edit).

$ws.Range("D2").Value = "'261.82"
$ws.Range("E2").Value = "'0.33%"
$ws.Range("D2:E2").ClearFormats()

$ws.Range("D3").Value = "'26.70"
$ws.Range("E3").Value = "'-1.71%"
$ws.Range("D3:E3").ClearFormats()

$ws.Range("D4").Value = "'4.704"
$ws.Range("E4").Value = "'0.32%"
$ws.Range("D4:E4").ClearFormats()

$ws.Range("D5").Value = "'0.06074"
$ws.Range("E5").Value = "'-0.87%"
$ws.Range("D5:E5").ClearFormats()

$ws.Range("D6").Value = "'6.701"
$ws.Range("E6").Value = "'0.66%"
$ws.Range("D6:E6").ClearFormats()

$ws.Range("D7").Value = "'0.8510"
$ws.Range("E7").Value = "'-0.28%"
$ws.Range("D7:E7").ClearFormats()

$ws.Range("D8").Value = "'0.9086"
$ws.Range("E8").Value = "'-1.33%"
$ws.Range("D8:E8").ClearFormats()

$ws.Range("D9").Value = "'0.1406"
$ws.Range("E9").Value = "'-0.16%"
$ws.Range("D9:E9").ClearFormats()

$ws.Range("D10").Value = "'0.05104"
$ws.Range("E10").Value = "'10.84%"
$ws.Range("D10:E10").ClearFormats()

$ws.Range("D11").Value = "'0.07090"
$ws.Range("E11").Value = "'0.03%"
$ws.Range("D11:E11").ClearFormats()

$ws.Range("D12").Value = "'0.03114"
$ws.Range("E12").Value = "'1.32%"
$ws.Range("D12:E12").ClearFormats()

$ws.Range("D13").Value = "'0.09046"
$ws.Range("E13").Value = "'-0.15%"
$ws.Range("D13:E13").ClearFormats()

$ws.Range("D14").Value = "'0.001530"
$ws.Range("E14").Value = "'0.09%"
$ws.Range("D14:E14").ClearFormats()

$ws.Range("D15").Value = "'0.0006168"
$ws.Range("E15").Value = "'1.61%"
$ws.Range("D15:E15").ClearFormats()

$ws.Range("D16").Value = "'0.006045"
$ws.Range("E16").Value = "'-0.34%"
$ws.Range("D16:E16").ClearFormats()

$ws.Range("D17").Value = "'3.450"
$ws.Range("E17").Value = "'-0.01%"
$ws.Range("D17:E17").ClearFormats()

$ws.Range("D18").Value = "'3.161"
$ws.Range("E18").Value = "'0.46%"
$ws.Range("D18:E18").ClearFormats()

$ws.Range("E19").Value = "'-0.80%"
$ws.Range("E19").ClearFormats()

$ws.Range("D22").Value = "'4.104"
$ws.Range("E22").Value = "'0.06%"
$ws.Range("D22:E22").ClearFormats()

$ws.Range("D23").Value = "'0.04242"
$ws.Range("E23").Value = "'-0.07%"
$ws.Range("D23:E23").ClearFormats()

$ws.Range("D24").Value = "'0.001180"
$ws.Range("E24").Value = "'-3.14%"
$ws.Range("D24:E24").ClearFormats()

$ws.Range("D25").Value = "'0.004060"
$ws.Range("E25").Value = "'6.90%"
$ws.Range("D25:E25").ClearFormats()

$ws.Range("E26").Value = "'0.03%"
$ws.Range("E26").ClearFormats()

$ws.Range("E27").Value = "'23.06%"
$ws.Range("E27").ClearFormats()

$ws.Range("D40").Value = "'0.03957"
$ws.Range("E40").Value = "'2.17%"
$ws.Range("D40:E40").ClearFormats()

$ws.Range("D41").Value = "'0.1111"
$ws.Range("E41").Value = "'-0.10%"
$ws.Range("D41:E41").ClearFormats()

$ws.Range("D42").Value = "'0.004188"
$ws.Range("E42").Value = "'2.64%"
$ws.Range("D42:E42").ClearFormats()

$ws.Range("D43").Value = "'0.002011"
$ws.Range("E43").Value = "'-9.28%"
$ws.Range("D43:E43").ClearFormats()

$ws.Range("D44").Value = "'0.01293"
$ws.Range("E44").Value = "'-20.86%"
$ws.Range("D44:E44").ClearFormats()

$ws.Range("D45").Value = "'0.00005109"
$ws.Range("E45").Value = "'-0.91%"
$ws.Range("D45:E45").ClearFormats()

$ws.Range("E46").Value = "'0.06%"
$ws.Range("E46").ClearFormats()

$ws.Range("D48").Value = "'0.2586"
$ws.Range("E48").Value = "'56.32%"
$ws.Range("D48:E48").ClearFormats()

$ws.Range("E49").Value = "'0.06%"
$ws.Range("E49").ClearFormats()

$ws.Range("E50").Value = "'0.06%"
$ws.Range("E50").ClearFormats()
